$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the title heading.
# -----------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# -----------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Fantastic Fireworks! Slot for
#    Free - Review") right before the closing "Prompt: ..." paragraph.
#    We do this by replacing the paragraph just before "Prompt" (the
#    "No progressive jackpot" bullet) with itself plus the new
#    paragraph - this inserts the new paragraph without disturbing the
#    "Prompt" paragraph that follows it.
# -----------------------------------------------------------------
$count = $d.Paragraphs.Count
$anchorPara = $d.Paragraphs.Item($count - 1)
$insertRange = $anchorPara.Range.Duplicate
$insertRange.Collapse(0)

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>No progressive jackpot</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fantastic Fireworks! Slot for Free - Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($xml)

# -----------------------------------------------------------------
# 3) Replace the old "Prompt: ..." text with the new meta-description
#    text, keeping the existing (italic) run formatting intact.
# -----------------------------------------------------------------
$oldText = "Prompt: Create a feature image that captures the excitement and explosive fun of the Fantastic Fireworks! slot game. The image should be in a cartoon style, and feature a happy Maya warrior wearing glasses. The warrior should be in a dynamic pose, perhaps with a lit sparkler in their hand, with colorful fireworks exploding behind them. The image should be eye-catching and convey the thrilling experience of playing Fantastic Fireworks! while also highlighting the unique elements of the game, such as the special fireworks and bonus features. Bonus points for incorporating the cityscape backdrop and the game's symbols, such as the fireworks characters."
$newText = "Find out if Fantastic Fireworks! is worth playing with our review. Play this high RTP game with explosive wins for free."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
